$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.614.27'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.873.31'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'246.81"
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = "'0.2909"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("D9").Value = "'0.06482"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = "'22.03"
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("D11").Value = "'0.07732"
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = "'0.7407"
$ws.Range("E12").Value = '  +3.60%  '
$ws.Range("D13").Value = "'96.49"
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").Value = '1.869.58'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = "'5.153"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = "'273.83"
$ws.Range("E16").Value = '  -0.88%  '
$ws.Range("D17").Value = '30.586.38'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = "'0.000007506"
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").Value = '2.116.76'
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = "'5.268"
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = "'6.197"
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").Value = "'9.257"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = "'163.55"
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").Value = "'18.81"
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D28").Value = "'1.918"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").Value = "'0.1001"
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").Value = "'1.355"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").Value = "'1.507"
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").Value = "'4.301"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").Value = "'4.116"
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D34").Value = "'0.04791"
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = "'1.119"
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").Value = "'0.6980"
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = "'0.9999"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = "'2.750"
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").Value = "'6.217"
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("D42").Value = "'73.55"
$ws.Range("E42").Value = '  +4.14%  '
$ws.Range("D43").Value = "'1.974"
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("D44").Value = "'0.4180"
$ws.Range("E44").Value = '  +1.45%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = "'0.8334"
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = "'102.58"
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = "'9.282"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").Value = "'35.37"
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = "'927.26"
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = "'6.954"
$ws.Range("E51").Value = '  -2.24%  '
